$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 2:
# F2 was "bookValue" -> becomes "total"
# G2 was "acb" -> becomes "bookValue"
$ws.Range("F2").Value = "total"
$ws.Range("G2").Value = "bookValue"

# Move the active selection to F3
$ws.Range("F3").Select()
